# Allow camp committee members to generate reports:
# Append a new suggestion row (row 6) to the suggestion list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 5
$ws.Cells.Item(6, 2).Value = 10
$ws.Cells.Item(6, 3).Value = "We should make poeple have fun"
$ws.Cells.Item(6, 4).Value = $true
$ws.Cells.Item(6, 5).Value = 7
$ws.Cells.Item(6, 6).Value = 2
